$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "145.30"), not numbers, in the
# source workbook. Force the cells whose new value looks like a plain
# decimal number to keep a Text number format so Excel does not silently
# convert them (and drop the significant trailing zero) when the value is
# assigned below.
$textPriceCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D17",
    "D19",
    "D20",
    "D22",
    "D23",
    "D26",
    "D28",
    "D29",
    "D30",
    "D32",
    "D36",
    "D37",
    "D40",
    "D41",
    "D43",
    "D46",
    "D48",
    "D50",
    "D51"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.124.31'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '2.928.43'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '592.02'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '145.30'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.505'
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").Value = '6.99'
$ws.Range("E9").Value = '  +5.29%  '
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '0.440'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '33.76'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '3.415.77'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '61.060.66'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '6.72'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").Value = '2.932.30'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").Value = '436.09'
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("D20").Value = '13.44'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '7.11'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").Value = '81.52'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("E24").Value = '  +3.16%  '
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '11.89'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  +2.71%  '
$ws.Range("D29").Value = '2.61'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("D32").Value = '26.67'
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '0.0₃0867'
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("D36").Value = '5.64'
$ws.Range("E36").Value = '  +1.17%  '
$ws.Range("D37").Value = '3.00'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = '8.59'
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").Value = '42.20'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("D43").Value = '376.97'
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '2.689.86'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = '133.40'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D48").Value = '24.08'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '2.00'
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").Value = '0.124'
$ws.Range("E51").Value = '  +0.41%  '
